# Data describing every new paragraph to append (Week 9 log entries),
# taken from the target diff: style, list level (-1 = not a list item,
# 0 = top-level bullet, 1 = sub-bullet), the paragraph text, and whether
# the run/paragraph-mark should be bold.
$paragraphs = @(
    @{ style = "Heading1"; ilvl = -1; text = "Week 9 (35 hours)"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "16/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added a copyright message to the main menu."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Fixed the rotation of the PlayerStart objects."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Re-prevented legged players from jumping on the log so it's in the build for testing sessions."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Tweaked some of the characters in the bones font to look more consistent in style."; bold = $true },
    @{ style = "ListParagraph"; ilvl = 1; text = "Repositioned the camera on the player so it doesn't play up when you have no legs."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Made the light only get brighter on the map in the last 30 seconds."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added a skylight and extensively tweaked the behaviour of the map brightening as the round ends."; bold = $true },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added some easter eggs."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "17/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Updated the new controls menu to show the controls on the controller image."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added some collision to the log tunnel so that legged players can walk down but not up it."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Fixed incorrect footstep sound/particle effect bug."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Made the player's own row on the scoreboard be highlighted."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "18/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Made stained glass windows destructible, fixed their replication, and added sound effects."; bold = $true },
    @{ style = "ListParagraph"; ilvl = 1; text = "Made it more satisfying to kill players (torso is now smashed)."; bold = $true },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added a bunch more graphical options to the options menu."; bold = $true },
    @{ style = "ListParagraph"; ilvl = 1; text = "Improved the sun rising sequence that happens near the end of the round."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added the UE4 intro movie to the game on startup."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "19/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Played the latest build of the game at a LAN party. Took notes on my friends' feedback and the bugs we encountered. "; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "20/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added a sound for dropping weapons."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Made the destructible torso show the correctly-coloured material"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Various usability tweaks to the menus."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added an option to toggle logging of gameplay statistics."; bold = $false },
    @{ style = "ListParagraph"; ilvl = 0; text = "21/09/2015"; bold = $false },
    @{ style = "ListParagraph"; ilvl = 1; text = "Added sound and particle effects to jumping."; bold = $true },
)

$d = $word.ActiveDocument

# The document currently ends with a trailing (empty, bookmark-only) paragraph.
# Every new paragraph must be inserted *before* it, i.e. right after the last
# real content paragraph ("Added sound effects to the bomb coffins ...").
$lastIndex = $d.Paragraphs.Count
$cursor = $d.Paragraphs.Item($lastIndex - 1)

foreach ($item in $paragraphs) {
    $newRange = $cursor.Range.InsertParagraphAfter()
    $cursor = $cursor.Next()
    $r = $cursor.Range

    $r.Text = $item.text
    $r.Style = $item.style

    if ($item.ilvl -ge 0) {
        # Word's ListLevelNumber is 1-based (1 = top level / ilvl 0, 2 = ilvl 1, ...)
        $r.ListFormat.ListLevelNumber = $item.ilvl + 1
    }

    if ($item.bold) {
        $r.Bold = 1
    }
}

Write-Output ("Inserted " + $paragraphs.Count + " paragraphs. New paragraph count: " + $d.Paragraphs.Count)
